{"js": "// Update the worksheet date and every two-digit x two-digit multiplication\n// answer cell to the new values from the latest generator run.\nconst replacements = [\n  [\"2025-11-21 Friday\", \"2025-11-22 Saturday\"],\n  [\"54\u00d730=1620\", \"55\u00d716=880\"],\n  [\"55\u00d791=5005\", \"16\u00d714=224\"],\n  [\"68\u00d798=6664\", \"19\u00d712=228\"],\n  [\"52\u00d781=4212\", \"55\u00d777=4235\"],\n  [\"43\u00d756=2408\", \"31\u00d742=1302\"],\n  [\"70\u00d791=6370\", \"73\u00d725=1825\"],\n  [\"41\u00d761=2501\", \"42\u00d749=2058\"],\n  [\"91\u00d775=6825\", \"37\u00d767=2479\"],\n  [\"16\u00d750=800\", \"86\u00d760=5160\"],\n  [\"48\u00d793=4464\", \"25\u00d766=1650\"],\n  [\"93\u00d782=7626\", \"80\u00d719=1520\"],\n  [\"89\u00d754=4806\", \"19\u00d777=1463\"],\n  [\"66\u00d774=4884\", \"21\u00d765=1365\"],\n  [\"45\u00d753=2385\", \"30\u00d760=1800\"],\n  [\"17\u00d756=952\", \"71\u00d795=6745\"],\n  [\"64\u00d764=4096\", \"14\u00d765=910\"],\n  [\"77\u00d771=5467\", \"66\u00d763=4158\"],\n  [\"87\u00d777=6699\", \"85\u00d779=6715\"],\n  [\"42\u00d715=630\", \"99\u00d755=5445\"],\n  [\"82\u00d745=3690\", \"49\u00d763=3087\"],\n  [\"12\u00d743=516\", \"77\u00d772=5544\"],\n  [\"18\u00d747=846\", \"47\u00d770=3290\"],\n  [\"42\u00d776=3192\", \"89\u00d718=1602\"],\n  [\"44\u00d780=3520\", \"27\u00d716=432\"],\n  [\"40\u00d716=640\", \"43\u00d730=1290\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit x two-digit multiplication\n# answer cell to the new values from the latest generator run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-21 Friday\", \"2025-11-22 Saturday\"),\n    @(\"54\u00d730=1620\", \"55\u00d716=880\"),\n    @(\"55\u00d791=5005\", \"16\u00d714=224\"),\n    @(\"68\u00d798=6664\", \"19\u00d712=228\"),\n    @(\"52\u00d781=4212\", \"55\u00d777=4235\"),\n    @(\"43\u00d756=2408\", \"31\u00d742=1302\"),\n    @(\"70\u00d791=6370\", \"73\u00d725=1825\"),\n    @(\"41\u00d761=2501\", \"42\u00d749=2058\"),\n    @(\"91\u00d775=6825\", \"37\u00d767=2479\"),\n    @(\"16\u00d750=800\",  \"86\u00d760=5160\"),\n    @(\"48\u00d793=4464\", \"25\u00d766=1650\"),\n    @(\"93\u00d782=7626\", \"80\u00d719=1520\"),\n    @(\"89\u00d754=4806\", \"19\u00d777=1463\"),\n    @(\"66\u00d774=4884\", \"21\u00d765=1365\"),\n    @(\"45\u00d753=2385\", \"30\u00d760=1800\"),\n    @(\"17\u00d756=952\",  \"71\u00d795=6745\"),\n    @(\"64\u00d764=4096\", \"14\u00d765=910\"),\n    @(\"77\u00d771=5467\", \"66\u00d763=4158\"),\n    @(\"87\u00d777=6699\", \"85\u00d779=6715\"),\n    @(\"42\u00d715=630\",  \"99\u00d755=5445\"),\n    @(\"82\u00d745=3690\", \"49\u00d763=3087\"),\n    @(\"12\u00d743=516\",  \"77\u00d772=5544\"),\n    @(\"18\u00d747=846\",  \"47\u00d770=3290\"),\n    @(\"42\u00d776=3192\", \"89\u00d718=1602\"),\n    @(\"44\u00d780=3520\", \"27\u00d716=432\"),\n    @(\"40\u00d716=640\",  \"43\u00d730=1290\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $old\n    $rng.Find.Replacement.Text = $new\n    $rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
